$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for the Fgf1-Cd44 LR-pair sheet (rows 2-10).
# Maps cell address -> new numeric value, taken from the updated TPM run.
$newValues = @{
    "G2" = 1.729797666666667
    "H2" = 5.189393000000001
    "I2" = 0.06436583050179444
    "J2" = 0.06436583050179444
    "M2" = 8.142376000000001
    "N2" = 24.427128
    "O2" = 0.1741313933276368
    "P2" = 0.1741313933276368
    "Q2" = 14.08466300592267
    "R2" = 126.761967053304
    "S2" = 0.01120811174796797
    "T2" = 0.01120811174796797
    "G3" = 1.729797666666667
    "H3" = 5.189393000000001
    "I3" = 0.06436583050179444
    "J3" = 0.06436583050179444
    "O3" = 0.5205382400466131
    "P3" = 0.5205382400466131
    "Q3" = 42.10387083366323
    "R3" = 378.9348375029691
    "S3" = 0.03350487612854269
    "T3" = 0.03350487612854269
    "G4" = 1.729797666666667
    "H4" = 5.189393000000001
    "I4" = 0.06436583050179444
    "J4" = 0.06436583050179444
    "O4" = 0.3053303666257501
    "P4" = 0.3053303666257501
    "Q4" = 24.696726059654
    "R4" = 222.270534536886
    "S4" = 0.01965284262528379
    "T4" = 0.01965284262528379
    "I5" = 0.2200595722726403
    "J5" = 0.2200595722726403
    "M5" = 8.142376000000001
    "N5" = 24.427128
    "O5" = 0.1741313933276368
    "P5" = 0.1741313933276368
    "Q5" = 48.15388681423467
    "R5" = 433.3849813281121
    "S5" = 0.03831927993491864
    "T5" = 0.03831927993491864
    "I6" = 0.2200595722726403
    "J6" = 0.2200595722726403
    "O6" = 0.5205382400466131
    "P6" = 0.5205382400466131
    "S6" = 0.1145494224562106
    "T6" = 0.1145494224562106
    "I7" = 0.2200595722726403
    "J7" = 0.2200595722726403
    "O7" = 0.3053303666257501
    "P7" = 0.3053303666257501
    "S7" = 0.06719086988151103
    "T7" = 0.06719086988151103
    "H8" = 57.69206699999999
    "I8" = 0.7155745972255653
    "J8" = 0.7155745972255653
    "M8" = 8.142376000000001
    "N8" = 24.427128
    "O8" = 0.1741313933276368
    "P8" = 0.1741313933276368
    "Q8" = 156.583500577064
    "R8" = 1409.251505193576
    "S8" = 0.1246040016447502
    "T8" = 0.1246040016447502
    "H9" = 57.69206699999999
    "I9" = 0.7155745972255653
    "J9" = 0.7155745972255653
    "O9" = 0.5205382400466131
    "P9" = 0.5205382400466131
    "Q9" = 468.0815920272456
    "R9" = 4212.734328245211
    "S9" = 0.3724839414618598
    "T9" = 0.3724839414618598
    "H10" = 57.69206699999999
    "I10" = 0.7155745972255653
    "J10" = 0.7155745972255653
    "O10" = 0.3053303666257501
    "P10" = 0.3053303666257501
    "S10" = 0.2184866541189553
    "T10" = 0.2184866541189553
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}
